$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Part_A")

# Pillar 6 / Q6 - natural hazard exposure assessment
$ws.Cells.Item(24, 1).Value = 6
$ws.Cells.Item(24, 2).Value = 2
$ws.Cells.Item(24, 3).Value = 1
$ws.Cells.Item(24, 4).Value = 'Has government assessed its exposure to natural hazards (e.g., earthquake, storm, flood etc.)?'
$ws.Cells.Item(24, 5).Value = 'yes/no'

$ws.Cells.Item(25, 1).Value = 6
$ws.Cells.Item(25, 2).Value = 2
$ws.Cells.Item(25, 3).Value = 2
$ws.Cells.Item(25, 4).Value = 'a. What type of assessment was undertaken (e.g. informal assessment or specialist assessment models such as probabilistic models)'
$ws.Cells.Item(25, 5).Value = 'text'

$ws.Cells.Item(26, 1).Value = 6
$ws.Cells.Item(26, 2).Value = 2
$ws.Cells.Item(26, 3).Value = 2
$ws.Cells.Item(26, 4).Value = 'b. If a specialist model (such as a probabilistic catastrophe model) was used, what assets (public, private etc.) and hazards were covered?'
$ws.Cells.Item(26, 5).Value = 'text'

$ws.Cells.Item(27, 1).Value = 6
$ws.Cells.Item(27, 2).Value = 2
$ws.Cells.Item(27, 3).Value = 2
$ws.Cells.Item(27, 4).Value = 'c. Other comments:'
$ws.Cells.Item(27, 5).Value = 'text'

$ws.Cells.Item(28, 1).Value = 7
$ws.Cells.Item(28, 2).Value = 2
$ws.Cells.Item(28, 3).Value = 1
$ws.Cells.Item(28, 4).Value = 'Are there any databases or registers which list or record the loss and damage to public assets caused by natural disasters?'
$ws.Cells.Item(28, 5).Value = 'yes/no'

$ws.Cells.Item(29, 1).Value = 7
$ws.Cells.Item(29, 2).Value = 2
$ws.Cells.Item(29, 3).Value = 2
$ws.Cells.Item(29, 4).Value = 'a. What information does it capture? (Please select all that apply)'
$ws.Cells.Item(29, 5).Value = 'yes/no'

$ws.Cells.Item(30, 1).Value = 7
$ws.Cells.Item(30, 2).Value = 2
$ws.Cells.Item(30, 3).Value = 2
$ws.Cells.Item(30, 4).Value = 'i. Public asset repair costs due to catastrophic events above a certain threshold'
$ws.Cells.Item(30, 5).Value = 'yes/no'

$ws.Cells.Item(31, 1).Value = 7
$ws.Cells.Item(31, 2).Value = 2
$ws.Cells.Item(31, 3).Value = 2
$ws.Cells.Item(31, 4).Value = 'ii. Public asset repair costs from all types of loss events including smaller loss events'
$ws.Cells.Item(31, 5).Value = 'yes/no'

$ws.Cells.Item(32, 1).Value = 7
$ws.Cells.Item(32, 2).Value = 2
$ws.Cells.Item(32, 3).Value = 2
$ws.Cells.Item(32, 4).Value = 'iii. The costs of public asset service interruption (e.g., duration of closure and costs associated with managing disruption) '
$ws.Cells.Item(32, 5).Value = 'yes/no'

$ws.Cells.Item(33, 1).Value = 7
$ws.Cells.Item(33, 2).Value = 2
$ws.Cells.Item(33, 3).Value = 2
$ws.Cells.Item(33, 4).Value = 'iv. Other:'
$ws.Cells.Item(33, 5).Value = 'yes/no'

$ws.Cells.Item(34, 1).Value = 8
$ws.Cells.Item(34, 2).Value = 2
$ws.Cells.Item(34, 3).Value = 1
$ws.Cells.Item(34, 4).Value = 'Has the government quantified the costs that it could incur from future natural disasters (i.e., the government’s contingent liabilities)?  '
$ws.Cells.Item(34, 5).Value = 'yes/no'

$ws.Cells.Item(35, 1).Value = 8
$ws.Cells.Item(35, 2).Value = 2
$ws.Cells.Item(35, 3).Value = 2
$ws.Cells.Item(35, 4).Value = 'If yes, please select all that apply:'

$ws.Cells.Item(36, 1).Value = 8
$ws.Cells.Item(36, 2).Value = 2
$ws.Cells.Item(36, 3).Value = 2
$ws.Cells.Item(36, 4).Value = 'a. The methodology to quantify the contingent liabilities is documented and can be replicated.'
$ws.Cells.Item(36, 5).Value = 'yes/no'

$ws.Cells.Item(37, 1).Value = 8
$ws.Cells.Item(37, 2).Value = 2
$ws.Cells.Item(37, 3).Value = 2
$ws.Cells.Item(37, 4).Value = 'b. The government accounts for the disaster related contingent liabilities from public assets in any fiscal forecasts, stress test or sensitivity analysis.'
$ws.Cells.Item(37, 5).Value = 'yes/no'

$ws.Cells.Item(38, 1).Value = 8
$ws.Cells.Item(38, 2).Value = 2
$ws.Cells.Item(38, 3).Value = 2
$ws.Cells.Item(38, 4).Value = 'c. The government publishes a fiscal risk statement that integrates disaster related contingent liabilities from public assets into broader fiscal planning.'
$ws.Cells.Item(38, 5).Value = 'yes/no'

$ws.Cells.Item(39, 1).Value = 8
$ws.Cells.Item(39, 2).Value = 2
$ws.Cells.Item(39, 3).Value = 2
$ws.Cells.Item(39, 4).Value = 'd. The government accounts for explicit government guarantees for asset damages incurred to public assets owned by subnational government, state owned-enterprises (SOEs) and/or operated through public-private partnerships (PPPs).'
$ws.Cells.Item(39, 5).Value = 'yes/no'

$ws.Cells.Item(40, 1).Value = 8
$ws.Cells.Item(40, 2).Value = 2
$ws.Cells.Item(40, 3).Value = 2
$ws.Cells.Item(40, 4).Value = 'e. Other:'
$ws.Cells.Item(40, 5).Value = 'yes/no'

$ws.Cells.Item(41, 1).Value = 9
$ws.Cells.Item(41, 2).Value = 3
$ws.Cells.Item(41, 3).Value = 1
$ws.Cells.Item(41, 4).Value = 'Does the government have risk management policies and/or measures in place to protect public assets prior to natural disasters occurring?'

$ws.Cells.Item(42, 1).Value = 9
$ws.Cells.Item(42, 2).Value = 3
$ws.Cells.Item(42, 3).Value = 2
$ws.Cells.Item(42, 4).Value = 'Please select all that apply:'

$ws.Cells.Item(43, 1).Value = 9
$ws.Cells.Item(43, 2).Value = 3
$ws.Cells.Item(43, 3).Value = 2
$ws.Cells.Item(43, 4).Value = 'a. The government has a plan or policy for making public assets more resilient against disaster and climate risks. '
$ws.Cells.Item(43, 5).Value = 'yes/no'

$ws.Cells.Item(44, 1).Value = 9
$ws.Cells.Item(44, 2).Value = 3
$ws.Cells.Item(44, 3).Value = 2
$ws.Cells.Item(44, 4).Value = 'b. The government has budget allocation processes in place to cover ongoing disaster risk management costs, and climate change adaptation spending for public assets.'
$ws.Cells.Item(44, 5).Value = 'yes/no'

$ws.Cells.Item(45, 1).Value = 9
$ws.Cells.Item(45, 2).Value = 3
$ws.Cells.Item(45, 3).Value = 2
$ws.Cells.Item(45, 4).Value = 'c. The government has a risk reduction investment plan for its public assets.'
$ws.Cells.Item(45, 5).Value = 'yes/no'

$ws.Cells.Item(46, 1).Value = 9
$ws.Cells.Item(46, 2).Value = 3
$ws.Cells.Item(46, 3).Value = 2
$ws.Cells.Item(46, 4).Value = 'd. The government has a system or process to track investments which reduce risks to public assets and make them more resilient.'
$ws.Cells.Item(46, 5).Value = 'yes/no'

$ws.Cells.Item(47, 1).Value = 9
$ws.Cells.Item(47, 2).Value = 3
$ws.Cells.Item(47, 3).Value = 2
$ws.Cells.Item(47, 4).Value = 'e. Other:'
$ws.Cells.Item(47, 5).Value = 'yes/no'

# Row height / wrap formatting consistent with existing multi-line rows
$ws.Rows.Item(26).RowHeight = 28.8
$ws.Rows.Item(34).RowHeight = 28.8
$ws.Rows.Item(37).RowHeight = 28.8
$ws.Rows.Item(38).RowHeight = 28.8
$ws.Rows.Item(39).RowHeight = 28.8
$ws.Rows.Item(41).RowHeight = 28.8
$ws.Rows.Item(44).RowHeight = 28.8

# Restore active selection to D3 as left by the author
$ws.Activate() | Out-Null
$ws.Range("D3").Select() | Out-Null

